$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; this shifts the existing rows 40-150 down to 41-151.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with its data.
$ws.Range("A40").Value = 7
$ws.Range("B40").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C40").Value = "Ñuble"
$ws.Range("D40").Value = 44525
$ws.Range("E40").Value = 16
$ws.Range("F40").Value = 100112017
$ws.Range("G40").Value = "Apio"
$ws.Range("H40").Value = "Americana (o)"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 60
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 8500
$ws.Range("M40").Value = 8250
$ws.Range("N40").Value = "`$/docena de matas"
$ws.Range("O40").Value = "Provincia del Elquí"
$ws.Range("P40").Value = 1375
$ws.Range("Q40").Value = 6
$ws.Range("R40").Value = "Hortaliza"
